# Regenerate the "K" column (G) values for the save_data sheet.
# The workbook's column G (header "K", previously derived from a
# "Strike#"-style count) is recalculated and rewritten with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G) for rows 3-27, replacing the old Strike#-based values.
$newK = @{
    3  = 1
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 4
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 0
    24 = 3
    25 = 2
    26 = 1
    27 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
